$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.33"
$ws.Range("E2").Value = "'0.68%"
$ws.Range("G2").Value = "'2"
$ws.Range("D3").Value = "'28.70"
$ws.Range("E3").Value = "'-2.86%"
$ws.Range("G3").Value = "'2"
$ws.Range("D4").Value = "'5.278"
$ws.Range("E4").Value = "'2.22%"
$ws.Range("G4").Value = "'2"
$ws.Range("D5").Value = "'0.05740"
$ws.Range("E5").Value = "'-0.17%"
$ws.Range("G5").Value = "'2"
$ws.Range("D6").Value = "'6.670"
$ws.Range("E6").Value = "'1.59%"
$ws.Range("G6").Value = "'2"
$ws.Range("D7").Value = "'3.232"
$ws.Range("E7").Value = "'3.83%"
$ws.Range("G7").Value = "'2"
$ws.Range("D8").Value = "'0.8610"
$ws.Range("E8").Value = "'0.41%"
$ws.Range("G8").Value = "'2"
$ws.Range("D9").Value = "'0.9292"
$ws.Range("E9").Value = "'8.47%"
$ws.Range("G9").Value = "'2"
$ws.Range("E10").Value = "'2.70%"
$ws.Range("G10").Value = "'2"
$ws.Range("D11").Value = "'0.07154"
$ws.Range("E11").Value = "'1.76%"
$ws.Range("G11").Value = "'2"
$ws.Range("D12").Value = "'0.03129"
$ws.Range("E12").Value = "'2.29%"
$ws.Range("G12").Value = "'2"
$ws.Range("D13").Value = "'0.09236"
$ws.Range("E13").Value = "'-1.25%"
$ws.Range("G13").Value = "'2"
$ws.Range("D14").Value = "'0.001524"
$ws.Range("E14").Value = "'-0.23%"
$ws.Range("G14").Value = "'2"
$ws.Range("D15").Value = "'0.0006055"
$ws.Range("E15").Value = "'0.74%"
$ws.Range("G15").Value = "'2"
$ws.Range("D16").Value = "'0.005876"
$ws.Range("E16").Value = "'-2.06%"
$ws.Range("G16").Value = "'2"
$ws.Range("D17").Value = "'3.506"
$ws.Range("E17").Value = "'0.33%"
$ws.Range("G17").Value = "'2"
$ws.Range("D18").Value = "'2.237"
$ws.Range("E18").Value = "'1.47%"
$ws.Range("G18").Value = "'2"
$ws.Range("E19").Value = "'-2.25%"
$ws.Range("G19").Value = "'2"
$ws.Range("D20").Value = "'0.03340"
$ws.Range("E20").Value = "'1.27%"
$ws.Range("G20").Value = "'2"
$ws.Range("E21").Value = "'2.79%"
$ws.Range("G21").Value = "'2"
$ws.Range("D22").Value = "'3.542"
$ws.Range("E22").Value = "'0.91%"
$ws.Range("G22").Value = "'2"
$ws.Range("D23").Value = "'0.04192"
$ws.Range("E23").Value = "'0.79%"
$ws.Range("G23").Value = "'2"
$ws.Range("D24").Value = "'0.1377"
$ws.Range("E24").Value = "'-0.21%"
$ws.Range("G24").Value = "'2"
$ws.Range("D25").Value = "'0.005035"
$ws.Range("E25").Value = "'21.86%"
$ws.Range("G25").Value = "'2"
$ws.Range("D26").Value = "'0.001222"
$ws.Range("E26").Value = "'-0.33%"
$ws.Range("G26").Value = "'2"
$ws.Range("E27").Value = "'-0.89%"
$ws.Range("G27").Value = "'2"
$ws.Range("D28").Value = "'0.0001936"
$ws.Range("E28").Value = "'33.62%"
$ws.Range("G28").Value = "'2"
$ws.Range("G29").Value = "'2"
$ws.Range("G30").Value = "'2"
$ws.Range("G31").Value = "'2"
$ws.Range("G32").Value = "'2"
$ws.Range("G33").Value = "'2"
$ws.Range("G34").Value = "'2"
$ws.Range("G35").Value = "'2"
$ws.Range("G36").Value = "'2"
$ws.Range("G37").Value = "'2"
$ws.Range("G38").Value = "'2"
$ws.Range("G39").Value = "'2"
$ws.Range("D40").Value = "'0.03838"
$ws.Range("E40").Value = "'3.02%"
$ws.Range("G40").Value = "'2"
$ws.Range("D41").Value = "'0.005683"
$ws.Range("E41").Value = "'62.41%"
$ws.Range("G41").Value = "'2"
$ws.Range("D42").Value = "'0.1079"
$ws.Range("E42").Value = "'0.90%"
$ws.Range("G42").Value = "'2"
$ws.Range("D43").Value = "'0.002198"
$ws.Range("E43").Value = "'-10.63%"
$ws.Range("G43").Value = "'2"
$ws.Range("D44").Value = "'0.009552"
$ws.Range("E44").Value = "'2.19%"
$ws.Range("G44").Value = "'2"
$ws.Range("D45").Value = "'0.00005277"
$ws.Range("E45").Value = "'0.06%"
$ws.Range("G45").Value = "'2"
$ws.Range("E46").Value = "'-0.10%"
$ws.Range("G46").Value = "'2"
$ws.Range("E47").Value = "'46.40%"
$ws.Range("G47").Value = "'2"
$ws.Range("D48").Value = "'0.002176"
$ws.Range("E48").Value = "'-11.16%"
$ws.Range("G48").Value = "'2"
$ws.Range("E49").Value = "'-0.10%"
$ws.Range("G49").Value = "'2"
$ws.Range("E50").Value = "'-0.10%"
$ws.Range("G50").Value = "'2"
$ws.Range("G51").Value = "'2"
